# Weekly update: insert this week's new price rows (Primera / Segunda)
# for "Acelga" at the top of the date block that starts at row 629,
# pushing all the previously-recorded weeks down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 629; everything that was on 629..660
# shifts down to 631..662 (Excel also carries the D-column date style
# down with it).
$ws.Rows("629:630").Insert()

# --- New row 629 : Calidad "Primera" ---
$ws.Cells.Item(629, 1).Value = 8
$ws.Cells.Item(629, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(629, 3).Value = "Coquimbo"
$ws.Cells.Item(629, 4).Value = 45075
$ws.Cells.Item(629, 5).Value = 4
$ws.Cells.Item(629, 6).Value = 100112009
$ws.Cells.Item(629, 7).Value = "Acelga"
$ws.Cells.Item(629, 8).Value = "Sin especificar"
$ws.Cells.Item(629, 9).Value = "Primera"
$ws.Cells.Item(629, 10).Value = 1600
$ws.Cells.Item(629, 11).Value = 450
$ws.Cells.Item(629, 12).Value = 500
$ws.Cells.Item(629, 13).Value = 475
$ws.Cells.Item(629, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(629, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(629, 16).Value = 238
$ws.Cells.Item(629, 17).Value = 2
$ws.Cells.Item(629, 18).Value = "Hortaliza"

# --- New row 630 : Calidad "Segunda" ---
$ws.Cells.Item(630, 1).Value = 8
$ws.Cells.Item(630, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(630, 3).Value = "Coquimbo"
$ws.Cells.Item(630, 4).Value = 45075
$ws.Cells.Item(630, 5).Value = 4
$ws.Cells.Item(630, 6).Value = 100112009
$ws.Cells.Item(630, 7).Value = "Acelga"
$ws.Cells.Item(630, 8).Value = "Sin especificar"
$ws.Cells.Item(630, 9).Value = "Segunda"
$ws.Cells.Item(630, 10).Value = 800
$ws.Cells.Item(630, 11).Value = 350
$ws.Cells.Item(630, 12).Value = 400
$ws.Cells.Item(630, 13).Value = 375
$ws.Cells.Item(630, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(630, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(630, 16).Value = 188
$ws.Cells.Item(630, 17).Value = 2
$ws.Cells.Item(630, 18).Value = "Hortaliza"
